# Apply the "add some question for each tag" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Rows 177-188: change tag in column A from "greet" to "ทักทาย"
#    (the questions in column B stay exactly the same)
for ($r = 177; $r -le 188; $r++) {
    $ws.Cells.Item($r, 1).Value = "ทักทาย"
}

# 2. Insert a new row of data at row 203: "สบายดี" / "บายดีป่าว"
$ws.Cells.Item(203, 1).Value = "สบายดี"
$ws.Cells.Item(203, 2).Value = "บายดีป่าว"

# 3. Append new "บอกลา" (goodbye) rows 223-227
$ws.Cells.Item(223, 1).Value = "บอกลา"
$ws.Cells.Item(223, 2).Value = "บาย"

$ws.Cells.Item(224, 1).Value = "บอกลา"
$ws.Cells.Item(224, 2).Value = "บ๊าย"

$ws.Cells.Item(225, 1).Value = "บอกลา"
$ws.Cells.Item(225, 2).Value = "บ้าย"

$ws.Cells.Item(226, 1).Value = "บอกลา"
$ws.Cells.Item(226, 2).Value = "ลาก่อน"

$ws.Cells.Item(227, 1).Value = "บอกลา"
$ws.Cells.Item(227, 2).Value = "ไว้เจอกันใหม่"
